$d = $word.ActiveDocument
$d.Content.Find.Execute("Jacqueline Leon", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Marshall Ye", 2)
